$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the table/relationship names to reflect the renamed entities
$ws.Range("D22").Value = "Supervisores_supervisa"
$ws.Range("A7").Value = "Supervisores_supervisa(Id_empleado,Id_supervisor,Nombre_supervisor)"
$ws.Range("D10").Value = "Empleados_tiene_proyectos"
$ws.Range("A3").Value = "Empleados_tiene_proyectos(Id_empleado,Id_proyecto)"

# Widen column A to fit the longer text
$ws.Columns("A").ColumnWidth = 66.3

# Update the active selection to A21
$ws.Range("A21").Select()
